$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45189 (2023-09-20) for every
# data row (rows 2 through 484). Bump it by one day to 45190 (2023-09-21).
$ws.Range("C2:C484").Value = 45190
